$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column A: reel positions re-numbered 0-9 (was 1-10)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# ---------------------------------------------------------------------------
# New notes about axis math (columns J/K on row 6) - entered first so the
# new shared-string table grows in the same order the original author typed
# things in.
# ---------------------------------------------------------------------------
$ws.Range("J6").Value = "if y = 180, sub all by 180"
$ws.Range("K6").Value = "use x - 180"

# ---------------------------------------------------------------------------
# Column C: map reel-position -> rotation digit (stored as text, like B)
# Entered in the same order as the original author (by reel value, not by
# row) so the shared-string table matches.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = "2"
$ws.Range("C7").Value = "1"
$ws.Range("C2").Value = "6"
$ws.Range("C9").Value = "3"
$ws.Range("C4").Value = "8"
$ws.Range("C5").Value = "9"
$ws.Range("C10").Value = "4"
$ws.Range("C3").Value = "7"
$ws.Range("C6").Value = "0"

# ---------------------------------------------------------------------------
# Column D: new raw rotation/angle values
# Columns E/F: blank "text" formatted placeholder cells (same style as col C)
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 16
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

$ws.Range("D3").Value = 17
$ws.Range("E3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"

$ws.Range("D4").Value = 18
$ws.Range("E4").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"

$ws.Range("D5").Value = 19
$ws.Range("E5").NumberFormat = "@"
$ws.Range("F5").NumberFormat = "@"

$ws.Range("D6").Value = 10
$ws.Range("E6").NumberFormat = "@"
$ws.Range("F6").NumberFormat = "@"

$ws.Range("D7").Value = 11
$ws.Range("E7").NumberFormat = "@"
$ws.Range("F7").NumberFormat = "@"

$ws.Range("D8").Value = 12
$ws.Range("E8").NumberFormat = "@"
$ws.Range("F8").NumberFormat = "@"

$ws.Range("D9").Value = 13
$ws.Range("E9").NumberFormat = "@"
$ws.Range("F9").NumberFormat = "@"

$ws.Range("D10").Value = 14
$ws.Range("E10").NumberFormat = "@"
$ws.Range("F10").NumberFormat = "@"

$ws.Range("D11").Value = 15
$ws.Range("E11").NumberFormat = "@"
$ws.Range("F11").NumberFormat = "@"

# ---------------------------------------------------------------------------
# New sample rotation values (columns M/N/O on row 8)
# ---------------------------------------------------------------------------
$ws.Range("M8").Value = -90
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 90

# ---------------------------------------------------------------------------
# Sheet view: stop scrolling to A4 / clear old selection, select B17 instead
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B17").Select()
